$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add new columns C..O ---
$headers = @(
    "Acceptance Rate",
    "Timeline",
    "Institution",
    "Program",
    "Degree Type",
    "Degree's Country of Origin",
    "Decision",
    "Notification",
    "Undergrad GPA",
    "GRE General",
    "GRE Verbal",
    "Analytical Writing",
    "Notes"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 3  # start at column C (3)
    $ws.Cells.Item(1, $col).Value = $headers[$i]
}

# Apply the existing header style (from A1) to the newly added header cells
# so no new style entries are introduced in styles.xml.
$ws.Range("A1").Copy()
$ws.Range("C1:O1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data rows 2-8, column A: new ID values ---
# These look numeric, so a plain .Value assignment would store them as
# numbers. Route them through a formula-that-returns-text, then convert
# to a static value via copy/paste-values so they land as shared strings
# (t="s") without touching cell styles (no quotePrefix / text numFmt).
$ids = @("985541", "985521", "985501", "985481", "985461", "985441", "985421")

for ($r = 0; $r -lt $ids.Length; $r++) {
    $row = $r + 2
    $cell = $ws.Cells.Item($row, 1)
    $cell.Formula = '="' + $ids[$r] + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163) | Out-Null  # xlPasteValues
}

$excel.CutCopyMode = $false
